# Auto-generated Excel COM-interop script
# Applies numeric cell updates (value changes, additions, and removals)
# to the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets per the target diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value2 = 4598.8335
$ws.Range("J32").Value2 = 4199.6665
$ws.Range("L32").Value2 = 4199.6665
$ws.Range("N32").Value2 = -4851.6665

$ws.Range("H69").Value2 = 9013
$ws.Range("I69").Value2 = 0
$ws.Range("J69").Value2 = 9013
$ws.Range("K69").Value2 = 0
$ws.Range("L69").Value2 = 27039
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value2 = -28787

$ws.Range("H72").Value2 = 9013
$ws.Range("I72").Value2 = 0
$ws.Range("J72").Value2 = 9013
$ws.Range("K72").Value2 = 0
$ws.Range("L72").Value2 = 81117
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value2 = -89853

$ws.Range("H107").Value2 = 347.1111
$ws.Range("I107").Value2 = 396.2857
$ws.Range("K107").Value2 = 396.2857
$ws.Range("M107").Value2 = 1523.7143

$ws.Range("H137").Value2 = 2279.1
$ws.Range("I137").Value2 = 1124.25
$ws.Range("J137").Value2 = 4011.375
$ws.Range("K137").Value2 = 3372.75
$ws.Range("L137").Value2 = 12034.125
$ws.Range("M137").Value2 = -822.75
$ws.Range("N137").Value2 = -17134.125

$ws.Range("H138").Value2 = 2834.0833
$ws.Range("J138").Value2 = 3092.673
$ws.Range("L138").Value2 = 9278.019
$ws.Range("N138").Value2 = -19558.019


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H8").Value2 = 967.2
$ws.Range("I8").Value2 = 776.6667
$ws.Range("J8").Value2 = 1253
$ws.Range("K8").Value2 = 776.6667
$ws.Range("L8").Value2 = 1253
$ws.Range("M8").Value2 = -632.6667
$ws.Range("N8").Value2 = -1541

$ws.Range("H43").Value2 = 49999
$ws.Range("J43").Value2 = 49999
$ws.Range("L43").Value2 = 49999
$ws.Range("N43").Value2 = -50625

$ws.Range("H74").Value2 = 3386.2856
$ws.Range("I74").Value2 = 1117.8572
$ws.Range("J74").Value2 = 5654.7144
$ws.Range("K74").Value2 = 1117.8572
$ws.Range("L74").Value2 = 5654.7144
$ws.Range("M74").Value2 = -243.8571999999999
$ws.Range("N74").Value2 = -7402.7144

$ws.Range("H77").Value2 = 3386.2856
$ws.Range("I77").Value2 = 1117.8572
$ws.Range("J77").Value2 = 5654.7144
$ws.Range("K77").Value2 = 5589.286
$ws.Range("L77").Value2 = 28273.572
$ws.Range("M77").Value2 = -1221.286
$ws.Range("N77").Value2 = -37009.572

$ws.Range("H94").Value2 = 1400
$ws.Range("I94").Value2 = 1400
$ws.Range("K94").Value2 = 1400
$ws.Range("M94").Value2 = -949


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value2 = 513.46155
$ws.Range("I22").Value2 = 541.25
$ws.Range("K22").Value2 = 541.25
$ws.Range("M22").Value2 = -368.25

$ws.Range("H86").Value2 = 4215.769
$ws.Range("J86").Value2 = 3942.125
$ws.Range("L86").Value2 = 3942.125
$ws.Range("N86").Value2 = -6188.125

$ws.Range("H89").Value2 = 4215.769
$ws.Range("J89").Value2 = 3942.125
$ws.Range("L89").Value2 = 19710.625
$ws.Range("N89").Value2 = -30942.625


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value2 = 4601.125
$ws.Range("I31").Value2 = 2233.25
$ws.Range("K31").Value2 = 2233.25
$ws.Range("M31").Value2 = -1938.25

$ws.Range("H34").Value2 = 4601.125
$ws.Range("I34").Value2 = 2233.25
$ws.Range("K34").Value2 = 2233.25
$ws.Range("M34").Value2 = -2031.25

$ws.Range("H43").Value2 = 31666.5
$ws.Range("J43").Value2 = 31666.5
$ws.Range("L43").Value2 = 31666.5
$ws.Range("N43").Value2 = -32034.5

$ws.Range("H62").Value2 = 43524.7
$ws.Range("J62").Value2 = 69907.836
$ws.Range("L62").Value2 = 69907.836
$ws.Range("N62").Value2 = -71155.836

$ws.Range("H65").Value2 = 43524.7
$ws.Range("J65").Value2 = 69907.836
$ws.Range("L65").Value2 = 349539.18
$ws.Range("N65").Value2 = -355779.18

$ws.Range("H99").Value2 = 11658.138
$ws.Range("J99").Value2 = 12897.685
$ws.Range("L99").Value2 = 12897.685
$ws.Range("N99").Value2 = -15893.685

$ws.Range("H101").Value2 = 31666.5
$ws.Range("J101").Value2 = 31666.5
$ws.Range("L101").Value2 = 31666.5
$ws.Range("N101").Value2 = -38156.5

$ws.Range("H126").Value2 = 11658.138
$ws.Range("J126").Value2 = 12897.685
$ws.Range("L126").Value2 = 38693.055
$ws.Range("N126").Value2 = -43633.055


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H46").Value2 = 350
$ws.Range("I46").Value2 = 0
$ws.Range("J46").Value2 = 350
$ws.Range("K46").Value2 = 0
$ws.Range("L46").Value2 = 1050
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value2 = -1232

$ws.Range("H98").Value2 = 745
$ws.Range("I98").Value2 = 0
$ws.Range("J98").Value2 = 745
$ws.Range("K98").Value2 = 0
$ws.Range("L98").Value2 = 2235
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value2 = -5231

$ws.Range("H107").Value2 = 1248.5834
$ws.Range("J107").Value2 = 1374.4445
$ws.Range("L107").Value2 = 4123.333500000001
$ws.Range("N107").Value2 = -7963.333500000001

$ws.Range("H138").Value2 = 3228.1428
$ws.Range("I138").Value2 = 1266.3334
$ws.Range("K138").Value2 = 3799.0002
$ws.Range("M138").Value2 = 1340.9998


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value2 = 1348
$ws.Range("J2").Value2 = 2072.6667
$ws.Range("L2").Value2 = 2072.6667
$ws.Range("N2").Value2 = -2298.6667

$ws.Range("H11").Value2 = 8401020
$ws.Range("J11").Value2 = 10000050
$ws.Range("L11").Value2 = 10000050
$ws.Range("N11").Value2 = -10000328


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value2 = 1040
$ws.Range("J22").Value2 = 2750
$ws.Range("L22").Value2 = 2750
$ws.Range("N22").Value2 = -3340

$ws.Range("H27").Value2 = 1040
$ws.Range("J27").Value2 = 2750
$ws.Range("L27").Value2 = 2750
$ws.Range("N27").Value2 = -2964

$ws.Range("H36").Value2 = 0
$ws.Range("J36").Value2 = 0
$ws.Range("L36").Value2 = 0
$ws.Range("N36").ClearContents()

$ws.Range("H55").Value2 = 417.5625
$ws.Range("I55").Value2 = 145.26666
$ws.Range("K55").Value2 = 145.26666
$ws.Range("M55").Value2 = 27.73334

$ws.Range("H61").Value2 = 4164.9165
$ws.Range("I61").Value2 = 4108.778
$ws.Range("J61").Value2 = 4333.3335
$ws.Range("K61").Value2 = 4108.778
$ws.Range("L61").Value2 = 4333.3335
$ws.Range("M61").Value2 = -3906.778
$ws.Range("N61").Value2 = -4737.3335

$ws.Range("H68").Value2 = 2062.8
$ws.Range("J68").Value2 = 2466.3333
$ws.Range("L68").Value2 = 2466.3333
$ws.Range("N68").Value2 = -3964.3333

$ws.Range("H71").Value2 = 2062.8
$ws.Range("J71").Value2 = 2466.3333
$ws.Range("L71").Value2 = 12331.6665
$ws.Range("N71").Value2 = -19819.6665

$ws.Range("H93").Value2 = 0
$ws.Range("I93").Value2 = 0
$ws.Range("J93").Value2 = 0
$ws.Range("K93").Value2 = 0
$ws.Range("L93").Value2 = 0
$ws.Range("M93").ClearContents()
$ws.Range("N93").ClearContents()

$ws.Range("H100").Value2 = 6996
$ws.Range("I100").Value2 = 994.6667
$ws.Range("K100").Value2 = 994.6667
$ws.Range("M100").Value2 = -453.6667

$ws.Range("H113").Value2 = 4164.9165
$ws.Range("I113").Value2 = 4108.778
$ws.Range("J113").Value2 = 4333.3335
$ws.Range("K113").Value2 = 4108.778
$ws.Range("L113").Value2 = 4333.3335
$ws.Range("M113").Value2 = -1938.778
$ws.Range("N113").Value2 = -8673.333500000001


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value2 = 8000
$ws.Range("J32").Value2 = 8000
$ws.Range("L32").Value2 = 8000
$ws.Range("N32").Value2 = -8634

$ws.Range("H100").Value2 = 1936.1818
$ws.Range("I100").Value2 = 2233
$ws.Range("J100").Value2 = 1580
$ws.Range("K100").Value2 = 4466
$ws.Range("L100").Value2 = 3160
$ws.Range("M100").Value2 = -3925
$ws.Range("N100").Value2 = -4242

$ws.Range("H113").Value2 = 1284.7142
$ws.Range("I113").Value2 = 1208.3334
$ws.Range("J113").Value2 = 1342
$ws.Range("K113").Value2 = 3625.0002
$ws.Range("L113").Value2 = 4026
$ws.Range("M113").Value2 = -1455.0002
$ws.Range("N113").Value2 = -8366

